$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.718.20"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.814.22"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.35"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.40"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.812.18"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.42"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  -2.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.45"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.456.67"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.866.70"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.799.06"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.89"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.77"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.702"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000155"
$ws.Range("E24").Value = "  +8.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.32"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.17"
$ws.Range("E26").Value = "  -3.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.96"
$ws.Range("E27").Value = "  -2.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.32"
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.12"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.18"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.18"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.768.70"
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.52"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.80"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.303"
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.54"
$ws.Range("E45").Value = "  +14.95%  "
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.84"
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.44"
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "147.15"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "392.55"
$ws.Range("E50").Value = "  -2.83%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.778.88"
$ws.Range("E51").Value = "  +3.81%  "
